$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-dashcam-videos/4006124/281474991205262/1743535585915/aNjCWRFZTU-camera-video-segment-driver-1743535590915.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEYZKPWSQ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T230005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIHorXjyZYrdISpbyW8k5BcuTHtrHz77OFA8UiD5OyhegAiEA4pJzrM5veaKsrQ8bY%2FdhiP05whCeQE%2FnNU7N4i8KjJoq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEy6aHpXxDuQoPWytSq6A0W%2FhZ%2Fr5qBg7Fkf5CT7zaJeEyBp4EW6xrbXBRaoqVjvPp4IDWaCUE2ZCzLPlqeo4HLs6J2ip8U31u0pyMweGbChCI6UNuaZ9ONxhcYT1ueAk7yz%2BGD9iGF%2FquPEp11B7vKBVIIboz%2BldmGKJPL6HJtspdNtptLiI0WxXhJNm5UeBE45Yudw9x82psYeWdk%2BUcR54UGBNTtQqxe1kKWqS1W8LbwBjmJRt78lKkpUB%2FnPqm4O67exFc9Sju2Ck9piMssR0910YiMnMg0gzKrE6yqMM2t6LPdH3%2BW%2F5abDubPL2BRzGleQ%2BrTgcG3G7KauhNnoSFOF6MiozwHF8uqSjikTinQOdFg6OujiHAswl8Ddgj0TilUsOvhAOTlQPMq%2B2E3Gm2eZ1UvdhiSQeZk%2FwE79k8PGeAiodJepi3lOmT2FcCsogfnSbSKONS4odRlQNcXcaAtY8c7K1VvzLgAeKwpodNfT0hseQitmsP9Ydn0Fw850df50JJzu%2Fyne7TMJYdjJ%2FLaASdyeB4bWPOh4i%2FFz5szJn5cjAWy2Aj0BraT%2B1ZYrYx51%2BGO62QGh34OiQHs6OPTirHoWQi8w8t22vwY6pQFlD127tMsDVCmBM3Zdcl1IU3WpC%2ByiZyOlM3qiGeddXIyUTozoshQPBkLmlkuJmo9e4G3l9PecAD%2BJ%2FUhP692bUErLW1w%2FiHv7rvdumVUfvYOtgxsV2cv1hWEgIezkTpmd5SJlk7XWiIHMR8rBAcOuioD17QwlCOdVBanbwMr8DmRrRG9pDCG03mkWi52naKARnuOOC2CCi081Luu6pezC2MaTISY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2007%3A00%3A05%20GMT&X-Amz-Signature=17c5bfbfb30356fc45957c55e5438b872b081afa9e9785308eccbea17043f6db'
$ws.Range("L3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535585915/3RxOLQEM5t-camera-video-segment-1743535590915.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEYZKPWSQ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T230005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIHorXjyZYrdISpbyW8k5BcuTHtrHz77OFA8UiD5OyhegAiEA4pJzrM5veaKsrQ8bY%2FdhiP05whCeQE%2FnNU7N4i8KjJoq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEy6aHpXxDuQoPWytSq6A0W%2FhZ%2Fr5qBg7Fkf5CT7zaJeEyBp4EW6xrbXBRaoqVjvPp4IDWaCUE2ZCzLPlqeo4HLs6J2ip8U31u0pyMweGbChCI6UNuaZ9ONxhcYT1ueAk7yz%2BGD9iGF%2FquPEp11B7vKBVIIboz%2BldmGKJPL6HJtspdNtptLiI0WxXhJNm5UeBE45Yudw9x82psYeWdk%2BUcR54UGBNTtQqxe1kKWqS1W8LbwBjmJRt78lKkpUB%2FnPqm4O67exFc9Sju2Ck9piMssR0910YiMnMg0gzKrE6yqMM2t6LPdH3%2BW%2F5abDubPL2BRzGleQ%2BrTgcG3G7KauhNnoSFOF6MiozwHF8uqSjikTinQOdFg6OujiHAswl8Ddgj0TilUsOvhAOTlQPMq%2B2E3Gm2eZ1UvdhiSQeZk%2FwE79k8PGeAiodJepi3lOmT2FcCsogfnSbSKONS4odRlQNcXcaAtY8c7K1VvzLgAeKwpodNfT0hseQitmsP9Ydn0Fw850df50JJzu%2Fyne7TMJYdjJ%2FLaASdyeB4bWPOh4i%2FFz5szJn5cjAWy2Aj0BraT%2B1ZYrYx51%2BGO62QGh34OiQHs6OPTirHoWQi8w8t22vwY6pQFlD127tMsDVCmBM3Zdcl1IU3WpC%2ByiZyOlM3qiGeddXIyUTozoshQPBkLmlkuJmo9e4G3l9PecAD%2BJ%2FUhP692bUErLW1w%2FiHv7rvdumVUfvYOtgxsV2cv1hWEgIezkTpmd5SJlk7XWiIHMR8rBAcOuioD17QwlCOdVBanbwMr8DmRrRG9pDCG03mkWi52naKARnuOOC2CCi081Luu6pezC2MaTISY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2007%3A00%3A05%20GMT&X-Amz-Signature=423c060eaec1989cac72aed606ef47a65627b4e69604e110b70dd83268b72ec3'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743535541388/Uf5dS3fS3z-camera-video-segment-driver-1743535543888.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEYZKPWSQ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T230005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIHorXjyZYrdISpbyW8k5BcuTHtrHz77OFA8UiD5OyhegAiEA4pJzrM5veaKsrQ8bY%2FdhiP05whCeQE%2FnNU7N4i8KjJoq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEy6aHpXxDuQoPWytSq6A0W%2FhZ%2Fr5qBg7Fkf5CT7zaJeEyBp4EW6xrbXBRaoqVjvPp4IDWaCUE2ZCzLPlqeo4HLs6J2ip8U31u0pyMweGbChCI6UNuaZ9ONxhcYT1ueAk7yz%2BGD9iGF%2FquPEp11B7vKBVIIboz%2BldmGKJPL6HJtspdNtptLiI0WxXhJNm5UeBE45Yudw9x82psYeWdk%2BUcR54UGBNTtQqxe1kKWqS1W8LbwBjmJRt78lKkpUB%2FnPqm4O67exFc9Sju2Ck9piMssR0910YiMnMg0gzKrE6yqMM2t6LPdH3%2BW%2F5abDubPL2BRzGleQ%2BrTgcG3G7KauhNnoSFOF6MiozwHF8uqSjikTinQOdFg6OujiHAswl8Ddgj0TilUsOvhAOTlQPMq%2B2E3Gm2eZ1UvdhiSQeZk%2FwE79k8PGeAiodJepi3lOmT2FcCsogfnSbSKONS4odRlQNcXcaAtY8c7K1VvzLgAeKwpodNfT0hseQitmsP9Ydn0Fw850df50JJzu%2Fyne7TMJYdjJ%2FLaASdyeB4bWPOh4i%2FFz5szJn5cjAWy2Aj0BraT%2B1ZYrYx51%2BGO62QGh34OiQHs6OPTirHoWQi8w8t22vwY6pQFlD127tMsDVCmBM3Zdcl1IU3WpC%2ByiZyOlM3qiGeddXIyUTozoshQPBkLmlkuJmo9e4G3l9PecAD%2BJ%2FUhP692bUErLW1w%2FiHv7rvdumVUfvYOtgxsV2cv1hWEgIezkTpmd5SJlk7XWiIHMR8rBAcOuioD17QwlCOdVBanbwMr8DmRrRG9pDCG03mkWi52naKARnuOOC2CCi081Luu6pezC2MaTISY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2007%3A00%3A05%20GMT&X-Amz-Signature=3b2a34b389330ec5b80e0f5c1a046429c881ff71a75e442f14718eb249e64d9a'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532925726/R4KzSHnOTc-camera-video-segment-driver-1743532928226.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEYZKPWSQ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T230005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIHorXjyZYrdISpbyW8k5BcuTHtrHz77OFA8UiD5OyhegAiEA4pJzrM5veaKsrQ8bY%2FdhiP05whCeQE%2FnNU7N4i8KjJoq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEy6aHpXxDuQoPWytSq6A0W%2FhZ%2Fr5qBg7Fkf5CT7zaJeEyBp4EW6xrbXBRaoqVjvPp4IDWaCUE2ZCzLPlqeo4HLs6J2ip8U31u0pyMweGbChCI6UNuaZ9ONxhcYT1ueAk7yz%2BGD9iGF%2FquPEp11B7vKBVIIboz%2BldmGKJPL6HJtspdNtptLiI0WxXhJNm5UeBE45Yudw9x82psYeWdk%2BUcR54UGBNTtQqxe1kKWqS1W8LbwBjmJRt78lKkpUB%2FnPqm4O67exFc9Sju2Ck9piMssR0910YiMnMg0gzKrE6yqMM2t6LPdH3%2BW%2F5abDubPL2BRzGleQ%2BrTgcG3G7KauhNnoSFOF6MiozwHF8uqSjikTinQOdFg6OujiHAswl8Ddgj0TilUsOvhAOTlQPMq%2B2E3Gm2eZ1UvdhiSQeZk%2FwE79k8PGeAiodJepi3lOmT2FcCsogfnSbSKONS4odRlQNcXcaAtY8c7K1VvzLgAeKwpodNfT0hseQitmsP9Ydn0Fw850df50JJzu%2Fyne7TMJYdjJ%2FLaASdyeB4bWPOh4i%2FFz5szJn5cjAWy2Aj0BraT%2B1ZYrYx51%2BGO62QGh34OiQHs6OPTirHoWQi8w8t22vwY6pQFlD127tMsDVCmBM3Zdcl1IU3WpC%2ByiZyOlM3qiGeddXIyUTozoshQPBkLmlkuJmo9e4G3l9PecAD%2BJ%2FUhP692bUErLW1w%2FiHv7rvdumVUfvYOtgxsV2cv1hWEgIezkTpmd5SJlk7XWiIHMR8rBAcOuioD17QwlCOdVBanbwMr8DmRrRG9pDCG03mkWi52naKARnuOOC2CCi081Luu6pezC2MaTISY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2007%3A00%3A05%20GMT&X-Amz-Signature=444e3b07399b13a621d8a3b9ccd60533d8d82685d079f82ce7308823473c3f70'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743532886131/8jJIZWs74D-camera-video-segment-driver-1743532888631.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSEYZKPWSQ%2F20250402%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250402T230005Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEHYaCXVzLXdlc3QtMiJHMEUCIHorXjyZYrdISpbyW8k5BcuTHtrHz77OFA8UiD5OyhegAiEA4pJzrM5veaKsrQ8bY%2FdhiP05whCeQE%2FnNU7N4i8KjJoq5gMI3%2F%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FARAEGgw3ODEyMDQ5NDIyNDQiDEy6aHpXxDuQoPWytSq6A0W%2FhZ%2Fr5qBg7Fkf5CT7zaJeEyBp4EW6xrbXBRaoqVjvPp4IDWaCUE2ZCzLPlqeo4HLs6J2ip8U31u0pyMweGbChCI6UNuaZ9ONxhcYT1ueAk7yz%2BGD9iGF%2FquPEp11B7vKBVIIboz%2BldmGKJPL6HJtspdNtptLiI0WxXhJNm5UeBE45Yudw9x82psYeWdk%2BUcR54UGBNTtQqxe1kKWqS1W8LbwBjmJRt78lKkpUB%2FnPqm4O67exFc9Sju2Ck9piMssR0910YiMnMg0gzKrE6yqMM2t6LPdH3%2BW%2F5abDubPL2BRzGleQ%2BrTgcG3G7KauhNnoSFOF6MiozwHF8uqSjikTinQOdFg6OujiHAswl8Ddgj0TilUsOvhAOTlQPMq%2B2E3Gm2eZ1UvdhiSQeZk%2FwE79k8PGeAiodJepi3lOmT2FcCsogfnSbSKONS4odRlQNcXcaAtY8c7K1VvzLgAeKwpodNfT0hseQitmsP9Ydn0Fw850df50JJzu%2Fyne7TMJYdjJ%2FLaASdyeB4bWPOh4i%2FFz5szJn5cjAWy2Aj0BraT%2B1ZYrYx51%2BGO62QGh34OiQHs6OPTirHoWQi8w8t22vwY6pQFlD127tMsDVCmBM3Zdcl1IU3WpC%2ByiZyOlM3qiGeddXIyUTozoshQPBkLmlkuJmo9e4G3l9PecAD%2BJ%2FUhP692bUErLW1w%2FiHv7rvdumVUfvYOtgxsV2cv1hWEgIezkTpmd5SJlk7XWiIHMR8rBAcOuioD17QwlCOdVBanbwMr8DmRrRG9pDCG03mkWi52naKARnuOOC2CCi081Luu6pezC2MaTISY%3D&X-Amz-SignedHeaders=host&response-expires=Thu%2C%2003%20Apr%202025%2007%3A00%3A05%20GMT&X-Amz-Signature=1a1633b4f08eb75033773407306d6191a4da3f827c0916798b3e7253e24525a7'
